$wb = $excel.ActiveWorkbook

# ---- Sheet: Summary ----
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.7209737827715356
$ws1.Range("C2").Value = 0.9244604316546763
$ws1.Range("D2").Value = 0.4812734082397004
$ws1.Range("E2").Value = 0.6330049261083743
$ws1.Range("F2").Value = 0.5323115161557581
$ws1.Range("G2").Value = 0.4903140592896977
$ws1.Range("H2").Value = 0.7209737827715357
$ws1.Range("I2").Value = 257
$ws1.Range("J2").Value = 21
$ws1.Range("K2").Value = 513
$ws1.Range("L2").Value = 277

# ---- Sheet: Classification Report ----
$ws2 = $wb.Worksheets.Item("Classification Report")

$ws2.Range("B2").Value = 0.649367088607595
$ws2.Range("C2").Value = 0.9606741573033708
$ws2.Range("D2").Value = 0.7749244712990937

$ws2.Range("B3").Value = 0.9244604316546763
$ws2.Range("C3").Value = 0.4812734082397004
$ws2.Range("D3").Value = 0.6330049261083743

$ws2.Range("B4").Value = 0.7209737827715356
$ws2.Range("C4").Value = 0.7209737827715356
$ws2.Range("D4").Value = 0.7209737827715356
$ws2.Range("E4").Value = 0.7209737827715356

$ws2.Range("B5").Value = 0.7869137601311356
$ws2.Range("C5").Value = 0.7209737827715356
$ws2.Range("D5").Value = 0.703964698703734

$ws2.Range("B6").Value = 0.7869137601311357
$ws2.Range("C6").Value = 0.7209737827715356
$ws2.Range("D6").Value = 0.703964698703734

# ---- Sheet: Confusion Matrix ----
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

$ws3.Range("B2").Value = 513
$ws3.Range("C2").Value = 21

$ws3.Range("B3").Value = 277
$ws3.Range("C3").Value = 257
